$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOE")

# Fix circular reference errors: SUM ranges previously included the
# formula's own row (and rows below it), causing self-referencing sums.
# Correct them to reference the actual data rows (E3:E31 / E3:E32 / etc.)

$ws.Range("D32").Formula = "=ROUND(SUM(E3:E31)*`$C`$30,0)"
$ws.Range("E32").Formula = "=ROUND(SUM(E3:E31)*0.25,0)"

$ws.Range("D33").Formula = "=ROUND(SUM(E3:E32)*0.20,0)"
$ws.Range("E33").Formula = "=ROUND(SUM(E3:E32)*0.20,0)"

$ws.Range("E34").Formula = "=SUM(E3:E33)"
$ws.Range("G34").Formula = "=TEXT(SUM(G3:G33),`"`$#,##0`")"
